# Generate Report for Handoff
# The localization-status report is regenerated: the rows for
# adb73576-fdd6-49ca-96f2-83f54a1e4446 and edb7b64d-b387-4334-ae49-3b80715bbcc2
# move ahead of 4cbda6e1-9396-404e-ae9a-df7f4d1ca222 (which is now ready for
# handoff) on every sheet (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: reorder the three rows and flip 4cbda6e1's status columns
# ---------------------------------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("A4").Value = "adb73576-fdd6-49ca-96f2-83f54a1e4446.md"
$ovw.Range("A5").Value = "edb7b64d-b387-4334-ae49-3b80715bbcc2.md"
$ovw.Range("A6").Value = "4cbda6e1-9396-404e-ae9a-df7f4d1ca222.md"
$ovw.Range("B6").Value = "Ready for handoff"
$ovw.Range("C6").Value = "Ready for handoff"

# ---------------------------------------------------------------------------
# zh-cn sheet: same reorder, plus handoff file/datetime details
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("A4").Value = "adb73576-fdd6-49ca-96f2-83f54a1e4446.md"
$zh.Range("C4").Value = "adb73576-fdd6-49ca-96f2-83f54a1e4446.124d8101d6ff29239d57603d2b23717faab74709.zh-cn.xlf"
$zh.Range("D4").Value = "2016-03-09 15:02:29"

$zh.Range("A5").Value = "edb7b64d-b387-4334-ae49-3b80715bbcc2.md"
$zh.Range("C5").Value = "edb7b64d-b387-4334-ae49-3b80715bbcc2.feee2006d607e5d0c4a04728cfa67b44ad4c2842.zh-cn.xlf"
$zh.Range("D5").Value = "2016-03-09 15:04:54"

$zh.Range("A6").Value = "4cbda6e1-9396-404e-ae9a-df7f4d1ca222.md"
$zh.Range("B6").Value = "Ready for handoff"
$zh.Range("C6").Value = "4cbda6e1-9396-404e-ae9a-df7f4d1ca222.fdae4a8ba869d2b46b3d3714d4cb9d8a2a763620.zh-cn.xlf"
$zh.Range("D6").Value = "2016-03-09 15:13:44"

# ---------------------------------------------------------------------------
# de-de sheet: same reorder, plus handoff file/datetime details
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("A4").Value = "adb73576-fdd6-49ca-96f2-83f54a1e4446.md"
$de.Range("C4").Value = "adb73576-fdd6-49ca-96f2-83f54a1e4446.124d8101d6ff29239d57603d2b23717faab74709.de-de.xlf"
$de.Range("D4").Value = "2016-03-09 15:02:38"

$de.Range("A5").Value = "edb7b64d-b387-4334-ae49-3b80715bbcc2.md"
$de.Range("C5").Value = "edb7b64d-b387-4334-ae49-3b80715bbcc2.feee2006d607e5d0c4a04728cfa67b44ad4c2842.de-de.xlf"
$de.Range("D5").Value = "2016-03-09 15:05:03"

$de.Range("A6").Value = "4cbda6e1-9396-404e-ae9a-df7f4d1ca222.md"
$de.Range("B6").Value = "Ready for handoff"
$de.Range("C6").Value = "4cbda6e1-9396-404e-ae9a-df7f4d1ca222.fdae4a8ba869d2b46b3d3714d4cb9d8a2a763620.de-de.xlf"
$de.Range("D6").Value = "2016-03-09 15:13:54"

Write-Host "Applied localization-status handoff report update"
